$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.009632597429711
$ws.Range("D2").Value = 1.01221784508477
$ws.Range("E2").Value = 1.012037756324119
$ws.Range("F2").Value = 1.012993391380036
$ws.Range("I2").Value = 1.022886706003748
$ws.Range("J2").Value = 1.014891502240957
$ws.Range("K2").Value = 1.015082666885817
$ws.Range("L2").Value = 1.014903120519523
$ws.Range("M2").Value = 1.015855879778966
$ws.Range("N2").Value = 1.009082513329069

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.011754289497215
$ws.Range("D3").Value = 1.014122975376818
$ws.Range("E3").Value = 1.013862720549221
$ws.Range("F3").Value = 1.016170763194075
$ws.Range("I3").Value = 1.02321240895198
$ws.Range("J3").Value = 1.016639393899762
$ws.Range("K3").Value = 1.016789190687002
$ws.Range("L3").Value = 1.016529661135266
$ws.Range("M3").Value = 1.018831285797483
$ws.Range("N3").Value = 1.00965935640418

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.013109463599595
$ws.Range("D4").Value = 1.015339716699707
$ws.Range("E4").Value = 1.015028681094754
$ws.Range("F4").Value = 1.018169970584031
$ws.Range("I4").Value = 1.023409264783684
$ws.Range("J4").Value = 1.017752782565885
$ws.Range("K4").Value = 1.017877077018535
$ws.Range("L4").Value = 1.017566863211157
$ws.Range("M4").Value = 1.020699877916072
$ws.Range("N4").Value = 1.010026734664452

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.013675036165484
$ws.Range("D5").Value = 1.015847488571206
$ws.Range("E5").Value = 1.01551536106557
$ws.Range("F5").Value = 1.018997101733472
$ws.Range("I5").Value = 1.023488730699065
$ws.Range("J5").Value = 1.018216721212888
$ws.Range("K5").Value = 1.018330592914295
$ws.Range("L5").Value = 1.017999322807996
$ws.Range("M5").Value = 1.021472105200651
$ws.Range("N5").Value = 1.010179802130534

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.013769758041059
$ws.Range("D6").Value = 1.015932528425117
$ws.Range("E6").Value = 1.015596874385632
$ws.Range("F6").Value = 1.019135205769795
$ws.Range("I6").Value = 1.023501881262126
$ws.Range("J6").Value = 1.01829437880677
$ws.Range("K6").Value = 1.018406517856745
$ws.Range("L6").Value = 1.018071726909556
$ws.Range("M6").Value = 1.021600990750387
$ws.Range("N6").Value = 1.010205422811736

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.013117036962089
$ws.Range("D7").Value = 1.015346516181553
$ws.Range("E7").Value = 1.015035197744933
$ws.Range("F7").Value = 1.018181074824411
$ws.Range("I7").Value = 1.023410339503545
$ws.Range("J7").Value = 1.017758997852334
$ws.Range("K7").Value = 1.017883151868847
$ws.Range("L7").Value = 1.017572655724746
$ws.Range("M7").Value = 1.020710248514615
$ws.Range("N7").Value = 1.010028785338179

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.010353357701123
$ws.Range("D8").Value = 1.012865059801865
$ws.Range("E8").Value = 1.012657649270008
$ws.Range("F8").Value = 1.01407911055771
$ws.Range("I8").Value = 1.022999677663768
$ws.Range("J8").Value = 1.015485911423566
$ws.Range("K8").Value = 1.015662830774272
$ws.Range("L8").Value = 1.015456029109697
$ws.Range("M8").Value = 1.01687332303834
$ws.Range("N8").Value = 1.009278695033773

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.00534326893047
$ws.Range("D9").Value = 1.008365747612006
$ws.Range("E9").Value = 1.008350044004354
$ws.Range("F9").Value = 1.0064039225024
$ws.Range("I9").Value = 1.022167916251895
$ws.Range("J9").Value = 1.011341378698207
$ws.Range("K9").Value = 1.011621187950972
$ws.Range("L9").Value = 1.011605538705742
$ws.Range("M9").Value = 1.009666168704558
$ws.Range("N9").Value = 1.007910545499057

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.00190221272335
$ws.Range("D10").Value = 1.005274947225769
$ws.Range("E10").Value = 1.005393201084254
$ws.Range("F10").Value = 1.000968143639332
$ws.Range("I10").Value = 1.021538102858346
$ws.Range("J10").Value = 1.008478641826419
$ws.Range("K10").Value = 1.008834038516612
$ws.Range("L10").Value = 1.008951840831332
$ws.Range("M10").Value = 1.004543760807116
$ws.Range("N10").Value = 1.006965195779484

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.000386605875622
$ws.Range("D11").Value = 1.003913478175164
$ws.Range("E11").Value = 1.004091285099985
$ws.Range("F11").Value = 0.9985342220951382
$ws.Range("I11").Value = 1.021246872656698
$ws.Range("J11").Value = 1.007213888122432
$ws.Range("K11").Value = 1.007603755121062
$ws.Range("L11").Value = 1.007780855370962
$ws.Range("M11").Value = 1.002246010773924
$ws.Range("N11").Value = 1.006547464531345

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 0.9998196450449423
$ws.Range("D12").Value = 1.003404157607465
$ws.Range("E12").Value = 1.003604325996574
$ws.Range("F12").Value = 0.9976176903666514
$ws.Range("I12").Value = 1.021135853606438
$ws.Range("J12").Value = 1.006740182554333
$ws.Range("K12").Value = 1.00714312236603
$ws.Range("L12").Value = 1.007342483561657
$ws.Range("M12").Value = 1.001380145537709
$ws.Range("N12").Value = 1.006390994567039

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 0.9999414434480693
$ws.Range("D13").Value = 1.003513574235866
$ws.Range("E13").Value = 1.003708934981477
$ws.Range("F13").Value = 0.9978148606065974
$ws.Range("I13").Value = 1.021159797265589
$ws.Range("J13").Value = 1.006841973700924
$ws.Range("K13").Value = 1.007242097050123
$ws.Range("L13").Value = 1.007436672429577
$ws.Range("M13").Value = 1.0015664434995
$ws.Range("N13").Value = 1.006424617769739

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.00033982315272
$ws.Range("D14").Value = 1.003871452038784
$ws.Range("E14").Value = 1.004051102407932
$ws.Range("F14").Value = 0.9984587180884922
$ws.Range("I14").Value = 1.021237754158497
$ws.Range("J14").Value = 1.007174812264493
$ws.Range("K14").Value = 1.007565754332918
$ws.Range("L14").Value = 1.007744689837546
$ws.Range("M14").Value = 1.002174693029252
$ws.Range("N14").Value = 1.006534557591354

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.000584743760855
$ws.Range("D15").Value = 1.004091469777026
$ws.Range("E15").Value = 1.004261472626118
$ws.Range("F15").Value = 0.9988537552598573
$ws.Range("I15").Value = 1.021285407310824
$ws.Range("J15").Value = 1.007379361400938
$ws.Range("K15").Value = 1.007764682459115
$ws.Range("L15").Value = 1.00793401310978
$ws.Range("M15").Value = 1.002547802844203
$ws.Range("N15").Value = 1.006602120673105

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.002002246474797
$ws.Range("D16").Value = 1.005364804745336
$ws.Range("E16").Value = 1.005479139573341
$ws.Range("F16").Value = 1.001127950495047
$ws.Range("I16").Value = 1.021557035599969
$ws.Range("J16").Value = 1.008562037466439
$ws.Range("K16").Value = 1.008915183775168
$ws.Range("L16").Value = 1.009029083203597
$ws.Range("M16").Value = 1.004694541716283
$ws.Range("N16").Value = 1.006992738655009

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.00288444993581
$ws.Range("D17").Value = 1.006157248349109
$ws.Range("E17").Value = 1.006237084287436
$ws.Range("F17").Value = 1.002532737790811
$ws.Range("I17").Value = 1.021722421616234
$ws.Range("J17").Value = 1.009297065412918
$ws.Range("K17").Value = 1.009630500465048
$ws.Range("L17").Value = 1.009710039653924
$ws.Range("M17").Value = 1.006019517289061
$ws.Range("N17").Value = 1.007235485784476

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.003396564214085
$ws.Range("D18").Value = 1.006617244951328
$ws.Range("E18").Value = 1.006677107110459
$ws.Range("F18").Value = 1.003344418621572
$ws.Range("I18").Value = 1.021817106266763
$ws.Range("J18").Value = 1.009723376068102
$ws.Range("K18").Value = 1.010045481269636
$ws.Range("L18").Value = 1.010105124094414
$ws.Range("M18").Value = 1.006784690588043
$ws.Range("N18").Value = 1.007376269940988

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.003570769074498
$ws.Range("D19").Value = 1.006773718917535
$ws.Range("E19").Value = 1.006826795541212
$ws.Range("F19").Value = 1.003619885611547
$ws.Range("I19").Value = 1.021849090772828
$ws.Range("J19").Value = 1.009868330782798
$ws.Range("K19").Value = 1.010186600970846
$ws.Range("L19").Value = 1.010239484080194
$ws.Range("M19").Value = 1.007044306974031
$ws.Range("N19").Value = 1.007424138305098

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.002790053364204
$ws.Range("D20").Value = 1.006072457466614
$ws.Range("E20").Value = 1.006155979356948
$ws.Range("F20").Value = 1.002382817586765
$ws.Range("I20").Value = 1.02170486203233
$ws.Range("J20").Value = 1.009218455020426
$ws.Range("K20").Value = 1.009553987537017
$ws.Range("L20").Value = 1.009637198162005
$ws.Range("M20").Value = 1.005878155341803
$ws.Range("N20").Value = 1.007209525018429

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.000222621845341
$ws.Range("D21").Value = 1.003766166736096
$ws.Range("E21").Value = 1.003950436704099
$ws.Range("F21").Value = 0.9982694656924104
$ws.Range("I21").Value = 1.021214876798272
$ws.Range("J21").Value = 1.007076908962195
$ws.Range("K21").Value = 1.007470547200005
$ws.Range("L21").Value = 1.007654081704021
$ws.Range("M21").Value = 1.001995923786297
$ws.Range("N21").Value = 1.006502219487144

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 0.9985851694371848
$ws.Range("D22").Value = 1.002295148733026
$ws.Range("E22").Value = 1.002544160657078
$ws.Range("F22").Value = 0.9956108923234179
$ws.Range("I22").Value = 1.020890319988964
$ws.Range("J22").Value = 1.005707680584361
$ws.Range("K22").Value = 1.006139412696714
$ws.Range("L22").Value = 1.006387388995422
$ws.Range("M22").Value = 0.9994831720325584
$ws.Range("N22").Value = 1.006049927652064

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 0.9994554661352381
$ws.Range("D23").Value = 1.003076997452071
$ws.Range("E23").Value = 1.003291553100484
$ws.Range("F23").Value = 0.9970272567167358
$ws.Range("I23").Value = 1.021063957864067
$ws.Range("J23").Value = 1.006435739728085
$ws.Range("K23").Value = 1.006847126986863
$ws.Range("L23").Value = 1.007060809416413
$ws.Range("M23").Value = 1.000822180609289
$ws.Range("N23").Value = 1.006290430709881

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.002832714703802
$ws.Range("D24").Value = 1.006110777675127
$ws.Range("E24").Value = 1.006192633583025
$ws.Range("F24").Value = 1.002450583863288
$ws.Range("I24").Value = 1.021712801955561
$ws.Range("J24").Value = 1.009253983138036
$ws.Range("K24").Value = 1.00958856738398
$ws.Range("L24").Value = 1.009670118596299
$ws.Range("M24").Value = 1.005942054371467
$ws.Range("N24").Value = 1.007221258059185

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.006655775098321
$ws.Range("D25").Value = 1.00954454195793
$ws.Range("E25").Value = 1.009478221313042
$ws.Range("F25").Value = 1.008442932397345
$ws.Range("I25").Value = 1.022396014383149
$ws.Range("J25").Value = 1.01242993445376
$ws.Range("K25").Value = 1.012681937196564
$ws.Range("L25").Value = 1.012615837051912
$ws.Range("M25").Value = 1.01158399435038
$ws.Range("N25").Value = 1.008269946482041
